$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: I1 = "I0", J1 = "IF" — copy H1's format (bold/border/center/top) onto them.
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2..37: I = 1 (constant), J = copy of column H's value.
for ($r = 2; $r -le 37; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
